$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 62 (pushing the
# existing rows 62-153 down to 63-154, which is why every row below it
# shows a one-row-later "previous" value in the diff).
$ws.Rows(62).Insert()

$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44477
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = 100112003
$ws.Range("G62").Value = "Ajo"
$ws.Range("H62").Value = "Chino"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 200
$ws.Range("K62").Value = 18000
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = 19000
$ws.Range("N62").Value = "`$/caja 10 kilos"
$ws.Range("O62").Value = "China"
$ws.Range("P62").Value = 1900
$ws.Range("Q62").Value = 10
$ws.Range("R62").Value = "Hortaliza"
